# Add cross-sectional Area columns (G, H) and a small summary block (J, K)
# to the discharge worksheet, mirroring the existing Q / Qtotal columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Area per segment (column G) -----------------------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Totals ----------------------------------------------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Summary block (J2/K2) --------------------------------------------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection / view state to mirror the saved workbook -------------------
$ws.Range("J2:K2").Select() | Out-Null
